$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value as literal text, avoiding Excels automatic
# numeric/percentage inference for strings that look like numbers, while
# leaving the cells style/number-format unchanged afterwards.
function Set-TextValue($cell, $text) {
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "42.000.64"
Set-TextValue $ws.Range("E2") "  -0.48%  "

Set-TextValue $ws.Range("D3") "2.215.35"
Set-TextValue $ws.Range("E3") "  -1.46%  "

Set-TextValue $ws.Range("E4") "  -0.06%  "

Set-TextValue $ws.Range("D5") "240.55"
Set-TextValue $ws.Range("E5") "  -2.56%  "

Set-TextValue $ws.Range("D6") "0.625"
Set-TextValue $ws.Range("E6") "  -0.81%  "

Set-TextValue $ws.Range("D7") "73.33"
Set-TextValue $ws.Range("E7") "  -1.44%  "

Set-TextValue $ws.Range("E8") "  +0.10%  "

Set-TextValue $ws.Range("D9") "0.607"
Set-TextValue $ws.Range("E9") "  -1.74%  "

Set-TextValue $ws.Range("D10") "42.68"
Set-TextValue $ws.Range("E10") "  +1.51%  "

Set-TextValue $ws.Range("D11") "0.0952"
Set-TextValue $ws.Range("E11") "  +0.65%  "

Set-TextValue $ws.Range("D12") "7.07"
Set-TextValue $ws.Range("E12") "  -1.74%  "

Set-TextValue $ws.Range("E13") "  -0.02%  "

Set-TextValue $ws.Range("D14") "2.548.85"
Set-TextValue $ws.Range("E14") "  -1.37%  "

Set-TextValue $ws.Range("D15") "14.27"
Set-TextValue $ws.Range("E15") "  -2.08%  "

Set-TextValue $ws.Range("D16") "0.837"
Set-TextValue $ws.Range("E16") "  -2.07%  "

Set-TextValue $ws.Range("D17") "2.214.85"
Set-TextValue $ws.Range("E17") "  -1.77%  "

Set-TextValue $ws.Range("D18") "41.838.04"
Set-TextValue $ws.Range("E18") "  -0.69%  "

Set-TextValue $ws.Range("E19") "  +7.57%  "

Set-TextValue $ws.Range("D20") "72.97"
Set-TextValue $ws.Range("E20") "  +0.73%  "

Set-TextValue $ws.Range("D21") "6.15"
Set-TextValue $ws.Range("E21") "  +0.35%  "

Set-TextValue $ws.Range("D22") "10.87"
Set-TextValue $ws.Range("E22") "  +22.00%  "

Set-TextValue $ws.Range("D23") "229.43"
Set-TextValue $ws.Range("E23") "  -0.89%  "

Set-TextValue $ws.Range("D24") "2.09"
Set-TextValue $ws.Range("E24") "  -6.50%  "

Set-TextValue $ws.Range("D25") "11.75"
Set-TextValue $ws.Range("E25") "  +2.82%  "

Set-TextValue $ws.Range("E26") "  +0.06%  "

Set-TextValue $ws.Range("E27") "  -0.22%  "

Set-TextValue $ws.Range("E28") "  -1.98%  "

Set-TextValue $ws.Range("E29") "  +0.92%  "

Set-TextValue $ws.Range("D30") "167.29"
Set-TextValue $ws.Range("E30") "  -2.39%  "

Set-TextValue $ws.Range("D31") "20.46"
Set-TextValue $ws.Range("E31") "  -1.04%  "

Set-TextValue $ws.Range("D32") "5.56"
Set-TextValue $ws.Range("E32") "  +6.93%  "

Set-TextValue $ws.Range("E33") "  -3.83%  "

Set-TextValue $ws.Range("D34") "29.61"
Set-TextValue $ws.Range("E34") "  -4.70%  "

Set-TextValue $ws.Range("D35") "0.125"
Set-TextValue $ws.Range("E35") "  -0.58%  "

Set-TextValue $ws.Range("E36") "  -9.84%  "

Set-TextValue $ws.Range("E37") "  -4.59%  "

Set-TextValue $ws.Range("E38") "  -5.17%  "

Set-TextValue $ws.Range("D39") "13.75"
Set-TextValue $ws.Range("E39") "  -0.95%  "

Set-TextValue $ws.Range("D40") "65.61"
Set-TextValue $ws.Range("E40") "  +4.36%  "

Set-TextValue $ws.Range("D41") "2.12"
Set-TextValue $ws.Range("E41") "  -2.77%  "

Set-TextValue $ws.Range("D42") "5.62"
Set-TextValue $ws.Range("E42") "  -2.89%  "

Set-TextValue $ws.Range("D43") "0.198"
Set-TextValue $ws.Range("E43") "  -3.42%  "

Set-TextValue $ws.Range("D44") "8.68"
Set-TextValue $ws.Range("E44") "  -0.10%  "

Set-TextValue $ws.Range("D45") "104.35"
Set-TextValue $ws.Range("E45") "  -2.22%  "

Set-TextValue $ws.Range("E46") "  -2.63%  "

Set-TextValue $ws.Range("D47") "2.38"
Set-TextValue $ws.Range("E47") "  +3.52%  "

Set-TextValue $ws.Range("D48") "1.12"
Set-TextValue $ws.Range("E48") "  -0.08%  "

Set-TextValue $ws.Range("E49") "  -0.66%  "

Set-TextValue $ws.Range("D50") "2.70"
Set-TextValue $ws.Range("E50") "  +0.20%  "

Set-TextValue $ws.Range("D51") "2.419.21"
Set-TextValue $ws.Range("E51") "  -1.66%  "
